# Auto-generated edit script applying the Coeurl_Profits.xlsx diff
# Updates per-row H..N metric cells across 8 job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# Cells that the diff deletes entirely (no replacement <v>) are cleared with ClearContents()
# rather than set to 0/blank, so the cell element itself is removed on save.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 4747.4863
$ws.Range("I138").Value = 4985
$ws.Range("J138").Value = 4740.8887
$ws.Range("K138").Value = 14955
$ws.Range("L138").Value = 14222.6661
$ws.Range("M138").Value = -9815
$ws.Range("N138").Value = -24502.6661

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1648.1625
$ws.Range("I32").Value = 1411.9324
$ws.Range("K32").Value = 1411.9324
$ws.Range("M32").Value = -1124.9324

$ws.Range("H61").Value = 3087.75
$ws.Range("I61").Value = 2942.8518
$ws.Range("K61").Value = 2942.8518
$ws.Range("M61").Value = -2730.8518

$ws.Range("H97").Value = 1623.3529
$ws.Range("I97").Value = 1623.3529
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 1623.3529
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -1127.3529
$ws.Range("N97").ClearContents()

$ws.Range("H102").Value = 5177.615
$ws.Range("I102").Value = 5976.125
$ws.Range("J102").Value = 3900
$ws.Range("K102").Value = 5976.125
$ws.Range("L102").Value = 3900
$ws.Range("M102").Value = -4354.125
$ws.Range("N102").Value = -7144

$ws.Range("H136").Value = 3087.75
$ws.Range("I136").Value = 2942.8518
$ws.Range("K136").Value = 8828.555399999999
$ws.Range("M136").Value = -6278.555399999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H36").Value = 1187.2
$ws.Range("I36").Value = 1187.2
$ws.Range("K36").Value = 1187.2
$ws.Range("M36").Value = -653.2

$ws.Range("H94").Value = 1295.8276
$ws.Range("I94").Value = 1245.84
$ws.Range("J94").Value = 1608.25
$ws.Range("K94").Value = 1245.84
$ws.Range("L94").Value = 1608.25
$ws.Range("M94").Value = -794.8399999999999
$ws.Range("N94").Value = -2510.25

$ws.Range("H99").Value = 3063.2188
$ws.Range("I99").Value = 1388.6154
$ws.Range("K99").Value = 1388.6154
$ws.Range("M99").Value = 109.3846000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4088799.5
$ws.Range("I31").Value = 6030481
$ws.Range("J31").Value = 11268
$ws.Range("K31").Value = 6030481
$ws.Range("L31").Value = 11268
$ws.Range("M31").Value = -6030186
$ws.Range("N31").Value = -11858

$ws.Range("H34").Value = 4088799.5
$ws.Range("I34").Value = 6030481
$ws.Range("J34").Value = 11268
$ws.Range("K34").Value = 6030481
$ws.Range("L34").Value = 11268
$ws.Range("M34").Value = -6030279
$ws.Range("N34").Value = -11672

$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()

$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()

$ws.Range("H105").Value = 952.7143
$ws.Range("I105").Value = 952.7143
$ws.Range("K105").Value = 952.7143
$ws.Range("M105").Value = 794.2857

$ws.Range("H132").Value = 4047.2903
$ws.Range("I132").Value = 3819.7144
$ws.Range("J132").Value = 6171.3335
$ws.Range("K132").Value = 11459.1432
$ws.Range("L132").Value = 18514.0005
$ws.Range("M132").Value = -8929.143199999999
$ws.Range("N132").Value = -23574.0005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 66669980
$ws.Range("I18").Value = 66669980
$ws.Range("K18").Value = 200009940
$ws.Range("M18").Value = -200009771

$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("M32").ClearContents()

$ws.Range("H56").Value = 11858.706
$ws.Range("I56").Value = 11858.706
$ws.Range("K56").Value = 11858.706
$ws.Range("M56").Value = -11328.706

$ws.Range("H68").Value = 1728429.1
$ws.Range("I68").Value = 4185.8335
$ws.Range("K68").Value = 12557.5005
$ws.Range("M68").Value = -11746.5005

$ws.Range("H71").Value = 1728429.1
$ws.Range("I71").Value = 4185.8335
$ws.Range("K71").Value = 37672.5015
$ws.Range("M71").Value = -33616.5015

$ws.Range("H107").Value = 17546144
$ws.Range("I107").Value = 55556300
$ws.Range("J107").Value = 2995.1538
$ws.Range("K107").Value = 166668900
$ws.Range("L107").Value = 8985.4614
$ws.Range("M107").Value = -166666980
$ws.Range("N107").Value = -12825.4614

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1625.5862
$ws.Range("I102").Value = 1575.68
$ws.Range("J102").Value = 1937.5
$ws.Range("K102").Value = 1575.68
$ws.Range("L102").Value = 1937.5
$ws.Range("M102").Value = 46.31999999999994
$ws.Range("N102").Value = -5181.5

$ws.Range("H132").Value = 2913.1396
$ws.Range("I132").Value = 2914.325
$ws.Range("J132").Value = 2897.3333
$ws.Range("K132").Value = 8742.974999999999
$ws.Range("L132").Value = 8691.999899999999
$ws.Range("M132").Value = -6212.974999999999
$ws.Range("N132").Value = -13751.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1015.1
$ws.Range("I22").Value = 1015.1
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 1015.1
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -720.1
$ws.Range("N22").ClearContents()

$ws.Range("H27").Value = 1015.1
$ws.Range("I27").Value = 1015.1
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 1015.1
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -908.1
$ws.Range("N27").ClearContents()

$ws.Range("H46").Value = 1963.4
$ws.Range("I46").Value = 1344.4286
$ws.Range("J46").Value = 5213
$ws.Range("K46").Value = 1344.4286
$ws.Range("L46").Value = 5213
$ws.Range("M46").Value = -1156.4286
$ws.Range("N46").Value = -5589

$ws.Range("H68").Value = 2999.8333
$ws.Range("I68").Value = 2999.8333
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 2999.8333
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -2250.8333
$ws.Range("N68").ClearContents()

$ws.Range("H71").Value = 2999.8333
$ws.Range("I71").Value = 2999.8333
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 14999.1665
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -11255.1665
$ws.Range("N71").ClearContents()

$ws.Range("H82").Value = 10955.4
$ws.Range("I82").Value = 15037
$ws.Range("J82").Value = 1431.6666
$ws.Range("K82").Value = 15037
$ws.Range("L82").Value = 1431.6666
$ws.Range("M82").Value = -14676
$ws.Range("N82").Value = -2153.6666

$ws.Range("H85").Value = 10955.4
$ws.Range("I85").Value = 15037
$ws.Range("J85").Value = 1431.6666
$ws.Range("K85").Value = 15037
$ws.Range("L85").Value = 1431.6666
$ws.Range("M85").Value = -13789
$ws.Range("N85").Value = -3927.6666

$ws.Range("H132").Value = 4047.6765
$ws.Range("I132").Value = 3626.05
$ws.Range("J132").Value = 4650
$ws.Range("K132").Value = 10878.15
$ws.Range("L132").Value = 13950
$ws.Range("M132").Value = -8348.150000000001
$ws.Range("N132").Value = -19010

$ws.Range("H136").Value = 3393.8845
$ws.Range("I136").Value = 3393.8845
$ws.Range("K136").Value = 10181.6535
$ws.Range("M136").Value = -7631.6535

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 10585.177
$ws.Range("J81").Value = 5071.4287
$ws.Range("L81").Value = 10142.8574
$ws.Range("N81").Value = -12264.8574

$ws.Range("H84").Value = 10585.177
$ws.Range("J84").Value = 5071.4287
$ws.Range("L84").Value = 50714.287
$ws.Range("N84").Value = -61322.287

$ws.Range("H132").Value = 1784.6296
$ws.Range("I132").Value = 1632.7084
$ws.Range("K132").Value = 4898.1252
$ws.Range("M132").Value = -2368.1252
